$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

# Heading / title text (two identical occurrences)
Replace-Text "Play Fortune of Camelot Free Online Slot Game" "Play Fortune of Camelot for Free"

# "What we like" bullet list
Replace-Text "Various bonus features, including Free Spins and stacked Wilds" "Multiple bonus features, including Free Spins and stacked Wilds"
Replace-Text "Impressive graphics that transport players to the magical world of Camelot" "User-friendly accessibility on mobile devices and PCs"
Replace-Text "Super Stake function doubles chance of getting Scatter symbols" "Impressive graphics that create an immersive experience"
Replace-Text "Playable on any mobile device or PC" "Turbo Spin and Super Stake functions for enhanced gameplay"

# "What we don't like" bullet list
Replace-Text "High volatility may not appeal to all players" "High volatility may result in less frequent wins"
Replace-Text "Limited betting range, minimum bet per spin is €0.20" "Limited number of paylines compared to some other slot games"

# Closing italic summary line
Replace-Text "Read the review of Fortune of Camelot and play this impressive online slot game for free. Enjoy various bonus features and immersive graphics." "Experience the magic of Camelot with Fortune of Camelot online slot game. Play now for free."
